# Weekly update: insert a new day's reading for Perejil (Mercado Mayorista
# Lo Valledor de Santiago) at row 822, pushing the existing historical rows
# (822-864) down by one to (823-865).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 822 - this shifts rows 822:864 down to 823:865
# and inherits formatting (incl. the date number format on column D) from
# the row above, matching how the rest of the sheet is laid out.
$ws.Rows("822:822").Insert()

# Populate the newly inserted row with this week's observation.
$ws.Range("A822").Value = 6
$ws.Range("B822").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C822").Value = "Metropolitana"
$ws.Range("D822").Value = 45267
$ws.Range("E822").Value = 13
$ws.Range("F822").Value = 100112044
$ws.Range("G822").Value = "Perejil"
$ws.Range("H822").Value = "Sin especificar"
$ws.Range("I822").Value = "Primera"
$ws.Range("J822").Value = 230
$ws.Range("K822").Value = 10000
$ws.Range("L822").Value = 12000
$ws.Range("M822").Value = 11130
$ws.Range("N822").Value = "`$/docena de atados"
$ws.Range("O822").Value = "Región Metropolitana"
$ws.Range("P822").Value = 3710
$ws.Range("Q822").Value = 3
$ws.Range("R822").Value = "Hortaliza"
